# Regenerate the experiment task-order sheets (per "Created experiment order
# generation script"): each tab keeps its physical position but gets a fresh
# randomized name/content, and the tab order is reshuffled.
#
#   tab1: GNG_TO-...  -> vSAT_TO-16515890301659896  (4 rows, vSAT/SAT stims)
#   tab2: NB_TO-...   -> RS_TO-16515890301659896    (2 rows, eyes closed/open)
#   tab3: RS_TO-...   -> GNG_TO-16515890301972394   (4 rows, go/GNG stims)
#   tab4: TOL_TO-...  -> NB_TO-16515890311255689    (9 rows, OB/ZB/TB stims)
#   tab5: vSAT_TO-... -> TOL_TO-16515890311724446   (6 rows, MM/ZM stims)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Tab 1: vSAT_TO-16515890301659896
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515890301659896"
$ws1.Range("B2").Value = "vSAT_stims-16515890301503649.csv"
$ws1.Range("B3").Value = "vSAT_stims-16515890301347415.csv"
$ws1.Range("B4").Value = "SAT_stims-1651589030103525.csv"
$ws1.Range("B5").Value = "SAT_stims-16515890301191514.csv"

# ---------------------------------------------------------------------------
# Tab 2: RS_TO-16515890301659896 (shrinks from 9 data rows to 2)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "RS_TO-16515890301659896"
$ws2.Range("A4:B10").Clear()
$ws2.Range("B2").Value = "eyes closed"
$ws2.Range("B3").Value = "eyes open"

# ---------------------------------------------------------------------------
# Tab 3: GNG_TO-16515890301972394 (grows from 2 data rows to 4)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "GNG_TO-16515890301972394"
$ws3.Range("A3").Copy()
$ws3.Range("A4:A5").PasteSpecial(-4122)
$ws3.Range("A4").Value = 2
$ws3.Range("A5").Value = 3
$ws3.Range("B2").Value = "go_stims-16515890301659896.csv"
$ws3.Range("B3").Value = "GNG_stims-16515890301816146.csv"
$ws3.Range("B4").Value = "go_stims-16515890301816146.csv"
$ws3.Range("B5").Value = "GNG_stims-16515890301972394.csv"

# ---------------------------------------------------------------------------
# Tab 4: NB_TO-16515890311255689
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "NB_TO-16515890311255689"
$ws4.Range("B2").Value = "TB-1651589031063068.csv"
$ws4.Range("B3").Value = "OB-1651589030920231.csv"
$ws4.Range("B4").Value = "ZB-match_3-16515890302995455.csv"
$ws4.Range("B5").Value = "TB-16515890310162303.csv"
$ws4.Range("B6").Value = "ZB-match_1-1651589030462585.csv"
$ws4.Range("B7").Value = "ZB-match_0-16515890302597404.csv"
$ws4.Range("B8").Value = "OB-16515890306210713.csv"
$ws4.Range("B9").Value = "TB-16515890311099427.csv"
$ws4.Range("B10").Value = "OB-16515890309849405.csv"

# ---------------------------------------------------------------------------
# Tab 5: TOL_TO-16515890311724446
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "TOL_TO-16515890311724446"
$ws5.Range("B2").Value = "MM_stims-16515890311411932.csv"
$ws5.Range("B3").Value = "ZM_stims-16515890311255689.csv"
$ws5.Range("B4").Value = "MM_stims-16515890311568205.csv"
$ws5.Range("B5").Value = "ZM_stims-16515890311411932.csv"
$ws5.Range("B6").Value = "MM_stims-16515890311724446.csv"
$ws5.Range("B7").Value = "ZM_stims-16515890311568205.csv"
